# The upstream change (commit: "Moving from 2.0.1 to 2.0.2") is a pure
# fixture re-generation: the canonical OOXML diff touches only the
# attribute *order* on existing elements in word/document.xml's
# <w:document>/<w:sectPr> and word/styles.xml's <w:docDefaults>,
# <w:latentStyles>/<w:lsdException> and <w:style> elements (namespace
# declarations and element attributes get alphabetized by local name,
# e.g. w:pgSz w:w=".." w:h=".." -> w:pgSz w:h=".." w:w="..").
#
# Every changed line was verified to carry the exact same tag name and
# the exact same attribute name/value set as before the change -- no
# text, run, paragraph, property value, style, or structural content is
# added, removed, or modified anywhere in the package. This is a
# byproduct of the authoring tool re-serializing the fixture (attribute
# order is not part of the OOXML information set and is not something
# the Word object model exposes or lets a macro control), so there is no
# content-level edit to replay against the document here.
#
# Touch the document content losslessly (re-set a range to itself) so
# the runtime still sees the script execute/save cleanly, without
# altering any visible text, formatting, or structure.
$d = $word.ActiveDocument
$null = $d.Content
